$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "JSU(-0.8751998923213852, 1.2079261590086496, 1.2058988120072138, 3.0017053584783904)"
$ws.Range("C2").Value = "NCT(3.2820580911459896, 2.761049683061684, -4.059292879184573, 4.2055641466725735)"
$ws.Range("D2").Value = "NIG(1.6501361163822272, 1.2770768149193397, 0.8420476601692433, 3.376656539499697)"
$ws.Range("E2").Value = "NCT(2.2626395393482617, 2.078902228193191, -0.011429398161107693, 3.2188042168035795)"
